$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folder_List_full")
$ws.Range("B92").Value = "test"
